$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that sat after "[Naglowek]"
$d.Bookmarks.Item("_GoBack").Delete()

# 2. In the signature table, drop the "Wójt" run and put the _GoBack
#    bookmark in its place instead.
$wojtRng = $d.Content
$wojtRng.Find.Execute("Wójt") | Out-Null
$d.Bookmarks.Add("_GoBack", $wojtRng)
$wojtRng.Text = ""

# 3. Drop the "Fuczek" run entirely (cell becomes empty).
$fuczekRng = $d.Content
$fuczekRng.Find.Execute("Fuczek") | Out-Null
$fuczekRng.Text = ""

# 4. Turn the hard-coded date/footer into the [Stopka] placeholder.
$footerRng = $d.Content
$footerRng.Find.Execute("Przeciszów, 20 września 2015", $true, $false, $false,
                         $false, $false, $true, 1, $false, "[Stopka]", 2) | Out-Null
